$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily columns being appended: GE=2025/03/14 ... GH=2025/03/17
$newDates = @("2025/03/14", "2025/03/15", "2025/03/16", "2025/03/17")
$newCols  = @(187, 188, 189, 190)   # GE, GF, GG, GH

# Per-row data for the 4 new columns (row, GE-value, GF-value, GG-value, GH-value).
# Cell shading in this sheet is value-driven: <125 -> yellow, 125-139.9 -> blue,
# >=140 -> no fill. Reference cells already carrying each look are reused below
# so new cells land on the exact same style (no new style entries created).
$rowData = @(
    @(2, 128.4, 154.7, 149.5, 160.5),
    @(3, 167.4, 264.5, 120.9, 166.7),
    @(4, 125.4, 144.8, 166, 121.9),
    @(5, 189.6, 166.9, 119.5, 363.2),
    @(6, 299.9, 198, 200.8, 142.2),
    @(7, 147.3, 162.9, 151.9, 123.3),
    @(8, 160.1, 150.6, 146.1, 207.9),
    @(9, 116.9, 139.2, 155.5, 139.3),
    @(10, 161.3, 158.9, 150.3, 185.7),
    @(11, 286.7, 151.5, 132.3, 141.6),
    @(12, 205.8, 166.3, 177.5, 152.2),
    @(13, 177.8, 134.1, 177.6, 183.1),
    @(14, 171.9, 164.3, 152.3, 115.1),
    @(15, 137.6, 129.3, 140.9, 188.5),
    @(16, 236.3, 150.3, 151, 126.4),
    @(17, 199.6, 155.4, 134.4, 158.7),
    @(18, 183.6, 133.8, 202.5, 155.9),
    @(19, 141.3, 134.8, 133.6, 144.9),
    @(20, 192.9, 178.4, 176.2, 140),
    @(21, 131.2, 144.8, 173.3, 176.3),
    @(22, 184.5, 133.2, 109.5, 218.6),
    @(23, 152.2, 143.7, 120.1, 156.2),
    @(24, 110.9, 134.9, 152.7, 149.5),
    @(25, 125.7, 175.3, 198, 255.3),
    @(26, 136.9, 135.7, 214.4, 121.6),
    @(27, 169.4, 165.8, 122.1, 143.3),
    @(28, 147.1, 148.9, 115.8, 131.8),
    @(29, 238.1, 174.9, 165.7, 169.4),
    @(30, 203.3, 139.8, 208.8, 133),
    @(31, 160.1, 149.8, 150.3, 145.2),
    @(32, 176.6, 131.2, 145.5, 163.8),
    @(33, 176.4, 152.5, 171.6, 136.7),
    @(34, 161.7, 144.1, 198.7, 159.3),
    @(35, 176.7, 169.9, 131.2, 186),
    @(36, 208.7, 142.4, 138.4, 151.4),
    @(37, 259.4, 122.1, 149.7, 127.7),
    @(38, 148.1, 187.9, 160.6, 118.3),
    @(39, 173.7, 133.8, 152.7, 141.2),
    @(40, 154.7, 169.6, 172.8, 148.5),
    @(41, 149.5, 140.9, 166.1, 121.5),
    @(42, 173.7, 149.4, 168.2, 187.3),
    @(43, 161.5, 142.8, 141, 130.2),
    @(44, 212.9, 139, 136.5, 141.4),
    @(45, 100.4, 147.8, 161.6, 198.7),
    @(46, 157.6, 169.4, 162.7, 192.7),
    @(47, 195.6, 128, 164.9, 155.3),
    @(48, 144, 200.2, 123.9, 173.4),
    @(49, 144.2, 213.7, 161.9, 148.9),
    @(50, 136.3, 161.5, 132.7, 192.1),
    @(51, 162.1, 139.8, 156.8, 141.4)
)

$styleRefLow    = "D2"   # < 125  -> yellow fill  (matches existing s=2 cells)
$styleRefMid    = "N2"   # 125-139.9 -> blue fill (matches existing s=3 cells)
$styleRefNormal = "B2"   # >= 140 -> default look (matches existing s=1 cells)

function Get-StyleRef($value) {
    if ($value -lt 125) { return $styleRefLow }
    elseif ($value -lt 140) { return $styleRefMid }
    else { return $styleRefNormal }
}

# --- 1. Give the 4 new columns the same width as all the others (width 12) ---
foreach ($col in $newCols) {
    $ws.Cells.Item(1, $col).ColumnWidth = 11.2
}

# --- 2. Header row: new date labels, stored as text (not auto-parsed as dates) ---
$ws.Range("GD1").Copy()
$ws.Cells.Item(1, $newCols[0]).PasteSpecial(-4122)
$ws.Cells.Item(1, $newCols[0]).NumberFormat = "@"
$ws.Cells.Item(1, $newCols[0]).Copy()
for ($i = 1; $i -lt $newCols.Length; $i++) {
    $ws.Cells.Item(1, $newCols[$i]).PasteSpecial(-4122)
}
for ($i = 0; $i -lt $newCols.Length; $i++) {
    $ws.Cells.Item(1, $newCols[$i]).Value = $newDates[$i]
}

# --- 3. Data rows: copy the matching style, then write the value ---
foreach ($entry in $rowData) {
    $r = $entry[0]
    for ($i = 0; $i -lt $newCols.Length; $i++) {
        $value = $entry[$i + 1]
        $ref = Get-StyleRef $value
        $ws.Range($ref).Copy()
        $target = $ws.Cells.Item($r, $newCols[$i])
        $target.PasteSpecial(-4122)
        $target.Value = $value
    }
}
